$p = $ppt.ActivePresentation

# Remove the final slide ("Thank you all, for your precious time / Happy Learning")
$lastIndex = $p.Slides.Count
$p.Slides.Item($lastIndex).Delete()
